# CDS_SPARSE_ENR.xlsx — "added new columns to sparse matrix"
#
# Enrollment_by_Race (sheet 2) gains six new sparse-matrix indicator
# columns (L:Q) — degree-seeking / non-degree-seeking / first-time /
# non-first-time / first-year / non-first-year — mirroring the layout
# already used on the General_Enrollment sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # General_Enrollment
$ws2 = $wb.Worksheets.Item(2)   # Enrollment_by_Race

# --- new header row (L1:Q1) ------------------------------------------------
$headers = New-Object 'object[,]' 1,6
$headers[0,0] = "degree-seeking"
$headers[0,1] = "non-degree-seeking"
$headers[0,2] = "first-time"
$headers[0,3] = "non-first-time"
$headers[0,4] = "first-year"
$headers[0,5] = "non-first-year"
$ws2.Range("L1:Q1").Value = $headers

# --- new data block (L2:Q21) ------------------------------------------------
# every enrollee is degree-seeking & first-time; first-year/non-first-year
# alternates with non-first-time, matching the source data rows.
$data = New-Object 'object[,]' 20,6
for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $odd = ($row % 2) -eq 1

    $data[$i,0] = 1   # degree-seeking
    $data[$i,1] = 0   # non-degree-seeking
    $data[$i,2] = 1   # first-time

    if ($odd) {
        $data[$i,3] = 1   # non-first-time
        $data[$i,4] = 1   # first-year
        $data[$i,5] = 1   # non-first-year
    } else {
        $data[$i,3] = 0   # non-first-time
        $data[$i,4] = 1   # first-year
        $data[$i,5] = 0   # non-first-year
    }
}
$ws2.Range("L2:Q21").Value = $data

# --- window/selection state -------------------------------------------------
# Enrollment_by_Race becomes the active tab; each sheet keeps its own
# last-used selection.
$ws1.Activate()
$ws1.Range("F3").Select() | Out-Null

$ws2.Activate()
$ws2.Range("G23").Select() | Out-Null

Write-Output "Added degree-seeking/first-time/first-year indicator columns (L:Q) to Enrollment_by_Race; set Enrollment_by_Race as the active tab."
